# Saldo.xlsx edit ("Add files via upload"):
# The account-balance row for 004467884 / ANA (balance 10100) is relocated
# further down the sheet, past the 004480970 / ALBERTO row, and its balance
# is corrected from 10100 to 100. Every row in between (005198093 .. the
# 004480970 row) shifts up by one row to close the gap left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$moveAccount = "004467884"
$stopAfterAccount = "004480970"

# Locate the row holding the account that needs to move, and the row of the
# account after which it should be reinserted.
$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$sourceRow = -1
$targetRow = -1
for ($r = 1; $r -le $lastRow; $r++) {
    $acct = $ws.Cells.Item($r, 1).Value2
    if ($acct -eq $moveAccount) {
        $sourceRow = $r
    }
    if ($acct -eq $stopAfterAccount) {
        $targetRow = $r
    }
}

if ($sourceRow -gt 0 -and $targetRow -gt 0 -and $targetRow -gt $sourceRow) {
    # Remember the row that is being displaced.
    $movedAccount = $ws.Cells.Item($sourceRow, 1).Value2
    $movedName = $ws.Cells.Item($sourceRow, 2).Value2

    # Shift every row below it, up through the target row, up by one.
    for ($r = $sourceRow + 1; $r -le $targetRow; $r++) {
        $acct = $ws.Cells.Item($r, 1).Value2
        $name = $ws.Cells.Item($r, 2).Value2
        $bal = $ws.Cells.Item($r, 3).Value2

        $dest = $r - 1
        $ws.Cells.Item($dest, 1).NumberFormat = "@"
        $ws.Cells.Item($dest, 1).Value = $acct
        $ws.Cells.Item($dest, 2).Value = $name
        $ws.Cells.Item($dest, 3).Value = $bal
    }

    # Re-insert the displaced row at the end of the block, with the
    # corrected balance.
    $ws.Cells.Item($targetRow, 1).NumberFormat = "@"
    $ws.Cells.Item($targetRow, 1).Value = $movedAccount
    $ws.Cells.Item($targetRow, 2).Value = $movedName
    $ws.Cells.Item($targetRow, 3).Value = 100
}
